$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.474.91'
$ws.Range('E2').Value = '  +2.53%  '

$ws.Range('D3').Value = '1.472.05'
$ws.Range('E3').Value = '  +4.06%  '

$ws.Range('E4').Value = '  +0.74%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.9628'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.82%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '275.90'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.08%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3652'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.04%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3062'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.48%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '39.82'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.03%  '

$ws.Range('E10').Value = '  +1.20%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.06607'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.51%  '

$ws.Range('E12').Value = '  +0.23%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '18.21'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +3.49%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.452'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.30%  '

$ws.Range('E15').Value = '  -0.31%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001028'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.94%  '

$ws.Range('D17').Value = '1.474.94'
$ws.Range('E17').Value = '  +4.24%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.05888'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +3.69%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.9713'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.95%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.02'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.67%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.449'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.49%  '

$ws.Range('E22').Value = '  -2.14%  '

$ws.Range('E23').Value = '  -0.29%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.248'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.53%  '

$ws.Range('D25').Value = '20.525.72'
$ws.Range('E25').Value = '  +2.77%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '141.78'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +6.63%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.126'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -6.54%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '17.18'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.12%  '

$ws.Range('D29').Value = '1.629.46'
$ws.Range('E29').Value = '  +3.47%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '113.34'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +3.29%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.882'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.92%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.949'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.71%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.8066'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.39%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.07868'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.27%  '

$ws.Range('B35').Value = 'TrustWalletToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.260'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +14.45%  '

$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.527'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +4.04%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.05729'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.50%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.746'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.74%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.9649'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.52%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.02038'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.33%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '7.641'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -5.01%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '10.40'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.39%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1876'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.40%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.5274'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.48%  '

$ws.Range('E45').Value = '  -0.94%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.01'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.68%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '116.73'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.05%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5167'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.28%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.765'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.03%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06443'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +4.23%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.9927'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.79%  '
